$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.618343710899353
$ws.Range("B1").Value = 2.591740369796753
$ws.Range("C1").Value = 2.867287397384644
$ws.Range("D1").Value = 3.2119460105896
$ws.Range("E1").Value = 2.946684598922729
